$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pdgfa"
$ws.Range("C2").Value = "Pdgfra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.668749
$ws.Range("H2").Value = 8.006247
$ws.Range("I2").Value = 0.08514925498774344
$ws.Range("J2").Value = 0.08514925498774345
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2781686666666667
$ws.Range("N2").Value = 0.834506
$ws.Range("O2").Value = 0.001228014730390642
$ws.Range("P2").Value = 0.001228014730390642
$ws.Range("Q2").Value = 0.7423623509980001
$ws.Range("R2").Value = 6.681261158982
$ws.Range("S2").Value = 0.0001045645394067378
$ws.Range("T2").Value = 0.0001045645394067378
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pdgfa"
$ws.Range("C3").Value = "Pdgfra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.668749
$ws.Range("H3").Value = 8.006247
$ws.Range("I3").Value = 0.08514925498774344
$ws.Range("J3").Value = 0.08514925498774345
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 225.778076
$ws.Range("N3").Value = 677.3342279999999
$ws.Range("O3").Value = 0.9967290940769435
$ws.Range("P3").Value = 0.9967290940769435
$ws.Range("Q3").Value = 602.5450145469239
$ws.Range("R3").Value = 5422.905130922316
$ws.Range("S3").Value = 0.08487073978526018
$ws.Range("T3").Value = 0.08487073978526019
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pdgfa"
$ws.Range("C4").Value = "Pdgfra"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.668749
$ws.Range("H4").Value = 8.006247
$ws.Range("I4").Value = 0.08514925498774344
$ws.Range("J4").Value = 0.08514925498774345
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4627536666666667
$ws.Range("N4").Value = 1.388261
$ws.Range("O4").Value = 0.002042891192665893
$ws.Range("P4").Value = 0.002042891192665893
$ws.Range("Q4").Value = 1.234973385163
$ws.Range("R4").Value = 11.114760466467
$ws.Range("S4").Value = 0.0001739506630765234
$ws.Range("T4").Value = 0.0001739506630765234
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Pdgfa"
$ws.Range("C5").Value = "Pdgfra"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.403995666666667
$ws.Range("H5").Value = 4.211987000000001
$ws.Range("I5").Value = 0.04479596433485759
$ws.Range("J5").Value = 0.04479596433485759
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2781686666666667
$ws.Range("N5").Value = 0.834506
$ws.Range("O5").Value = 0.001228014730390642
$ws.Range("P5").Value = 0.001228014730390642
$ws.Range("Q5").Value = 0.3905476026024445
$ws.Range("R5").Value = 3.514928423422
$ws.Range("S5").Value = [double]"5.501010406525898E-05"
$ws.Range("T5").Value = [double]"5.501010406525898E-05"
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Pdgfa"
$ws.Range("C6").Value = "Pdgfra"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.403995666666667
$ws.Range("H6").Value = 4.211987000000001
$ws.Range("I6").Value = 0.04479596433485759
$ws.Range("J6").Value = 0.04479596433485759
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 225.778076
$ws.Range("N6").Value = 677.3342279999999
$ws.Range("O6").Value = 0.9967290940769435
$ws.Range("P6").Value = 0.9967290940769435
$ws.Range("Q6").Value = 316.9914403323373
$ws.Range("R6").Value = 2852.922962991036
$ws.Range("S6").Value = 0.04464944094978568
$ws.Range("T6").Value = 0.04464944094978568
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Pdgfa"
$ws.Range("C7").Value = "Pdgfra"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.403995666666667
$ws.Range("H7").Value = 4.211987000000001
$ws.Range("I7").Value = 0.04479596433485759
$ws.Range("J7").Value = 0.04479596433485759
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4627536666666667
$ws.Range("N7").Value = 1.388261
$ws.Range("O7").Value = 0.002042891192665893
$ws.Range("P7").Value = 0.002042891192665893
$ws.Range("Q7").Value = 0.6497041427341113
$ws.Range("R7").Value = 5.847337284607001
$ws.Range("S7").Value = [double]"9.151328100665603E-05"
$ws.Range("T7").Value = [double]"9.151328100665602E-05"
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Pdgfa"
$ws.Range("C8").Value = "Pdgfra"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 27.26926766666667
$ws.Range("H8").Value = 81.80780300000001
$ws.Range("I8").Value = 0.8700547806773989
$ws.Range("J8").Value = 0.8700547806773991
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2781686666666667
$ws.Range("N8").Value = 0.834506
$ws.Range("O8").Value = 0.001228014730390642
$ws.Range("P8").Value = 0.001228014730390642
$ws.Range("Q8").Value = 7.585455827813112
$ws.Range("R8").Value = 68.269102450318
$ws.Range("S8").Value = 0.001068440086918646
$ws.Range("T8").Value = 0.001068440086918646
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Pdgfa"
$ws.Range("C9").Value = "Pdgfra"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 27.26926766666667
$ws.Range("H9").Value = 81.80780300000001
$ws.Range("I9").Value = 0.8700547806773989
$ws.Range("J9").Value = 0.8700547806773991
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 225.778076
$ws.Range("N9").Value = 677.3342279999999
$ws.Range("O9").Value = 0.9967290940769435
$ws.Range("P9").Value = 0.9967290940769435
$ws.Range("Q9").Value = 6156.802787709009
$ws.Range("R9").Value = 55411.22508938109
$ws.Range("S9").Value = 0.8672089133418976
$ws.Range("T9").Value = 0.8672089133418978
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Pdgfa"
$ws.Range("C10").Value = "Pdgfra"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 27.26926766666667
$ws.Range("H10").Value = 81.80780300000001
$ws.Range("I10").Value = 0.8700547806773989
$ws.Range("J10").Value = 0.8700547806773991
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4627536666666667
$ws.Range("N10").Value = 1.388261
$ws.Range("O10").Value = 0.002042891192665893
$ws.Range("P10").Value = 0.002042891192665893
$ws.Range("Q10").Value = 12.61895360006478
$ws.Range("R10").Value = 113.570582400583
$ws.Range("S10").Value = 0.001777427248582713
$ws.Range("T10").Value = 0.001777427248582713
